$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update the three parameter values that changed on the "Main" sheet.
# C7  (OFtarget)       : 5                  -> 10.199999999999999
# C17 (CdOx)            : 0.26369999999999999 -> 0.29060000000000002
# C18 (CdFuel)           : 0.1797             -> 0.13719999999999999
$ws.Range("C7").Value = 10.199999999999999
$ws.Range("C17").Value = 0.29060000000000002
$ws.Range("C18").Value = 0.13719999999999999

# Move the active selection to E14, matching the saved view state.
$ws.Range("E14").Select()

$excel.Calculate()
